$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the TYPE cell for the ETH-USD row from "Custom" to "Technical"
$ws.Range("B2").Value = "Technical"

# Move the active selection from F2 to B2
$ws.Range("B2").Select()
